$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 166; this shifts existing rows 166-223 down to 167-224
$ws.Rows.Item(166).Insert()

# Populate the newly inserted row 166 with the new data record
$ws.Cells.Item(166, 1).Value = 7
$ws.Cells.Item(166, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(166, 3).Value = "Ñuble"
$ws.Cells.Item(166, 4).Value = 44588
$ws.Cells.Item(166, 5).Value = 16
$ws.Cells.Item(166, 6).Value = 100112008
$ws.Cells.Item(166, 7).Value = "Coliflor"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 160
$ws.Cells.Item(166, 11).Value = 800
$ws.Cells.Item(166, 12).Value = 850
$ws.Cells.Item(166, 13).Value = 825
$ws.Cells.Item(166, 14).Value = '$/unidad'
$ws.Cells.Item(166, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(166, 16).Value = 825
$ws.Cells.Item(166, 17).Value = 1
$ws.Cells.Item(166, 18).Value = "Hortaliza"
